$wb = $excel.ActiveWorkbook

# "books" sheet: update the ISBN column header text
$booksWs = $wb.Worksheets.Item("books")
$booksWs.Range("C1").Value = "ISBNs"

# Update the active cell selection on the "books" sheet
$booksWs.Activate()
$booksWs.Range("C23").Select()
